$wb = $excel.ActiveWorkbook

# "Forecast Comparison" sheet — update yhat_upper (column D) forecasts,
# and a handful of rounded-to-int "Prophet Forecast" values (column B)
# that shifted with them.
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Range("D2").Value = -12.0649011511458
$ws1.Range("D3").Value = -58.60255706277554
$ws1.Range("D4").Value = -82.17838698355609
$ws1.Range("D5").Value = -30.02321168508124
$ws1.Range("D6").Value = 85.62308419655966
$ws1.Range("D7").Value = 186.0589138553771
$ws1.Range("D8").Value = 179.9279473612815
$ws1.Range("D9").Value = 69.74082401413624
$ws1.Range("D10").Value = -35.82501465321602
$ws1.Range("D11").Value = -21.41693286155797

$ws1.Range("B12").Value = 96
$ws1.Range("D12").Value = 103.1766388522287

$ws1.Range("D13").Value = 207.2844332072682
$ws1.Range("D14").Value = 175.877213565638

$ws1.Range("B15").Value = 32
$ws1.Range("D15").Value = 39.32749247658997

$ws1.Range("D16").Value = -50.14715971954683

$ws1.Range("B17").Value = 7
$ws1.Range("D17").Value = 14.84852449256483

$ws1.Range("B18").Value = 160
$ws1.Range("D18").Value = 167.6023283419961

$ws1.Range("D19").Value = 244.121833014396
$ws1.Range("D20").Value = 160.1105965927824
$ws1.Range("D21").Value = 2.525499866017072

# "Summary" sheet — total 16-week forecast ticked up by one unit.
# (kept as text, matching the existing column's string-typed cells)
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "994"
